# Update Register Account Test Case Script
# Adds TS_03 ("Not fill lastname inputbox") and TS_04 ("fill all fields but
# skip the privacy-policy checkbox") scenarios to the QAFOX register-account
# test-case sheet, widens column B, tweaks the active selection/view, and
# trims one trailing blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Column B a little wider to fit the new scenario titles.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 32.6

# ---------------------------------------------------------------------------
# 2) TS_03 block (rows 20-33) - "Validate Register Account Remains one
#    mandatory field should be blank" - identical flow to TS_02 except the
#    lastname inputbox is deliberately left blank.
# ---------------------------------------------------------------------------
$ws.Cells.Item(20,1).Value = 'TS_03'
$ws.Cells.Item(20,2).Value = 'Validate Register Account Remains one mandatory field should be blank'
$ws.Cells.Item(20,3).Value = 'Pre-requisite: TS_01'
$ws.Cells.Item(20,3).Font.Bold = $true
$ws.Cells.Item(20,3).Borders.LineStyle = 1
$ws.Cells.Item(20,3).HorizontalAlignment = -4131
$ws.Cells.Item(20,3).VerticalAlignment = -4160
$ws.Cells.Item(20,3).WrapText = $true

$ws.Cells.Item(21,3).Value = 'Click firstname inputbox'
$ws.Cells.Item(21,4).Value = 'firstname inputbox should be clicked successfully'

$ws.Cells.Item(22,3).Value = 'Enter valid firstname'
$ws.Cells.Item(22,4).Value = 'Firstname should be entered successfully (Entered firstname should be displayed in firstname inputbox)'

$ws.Cells.Item(23,3).Value = 'Not fill lastname inputbox'
$ws.Cells.Item(23,4).Value = 'Lastname  inputbox should be displayed empty'

$ws.Cells.Item(24,3).Value = 'Click email inputbox'
$ws.Cells.Item(24,4).Value = 'Email inputbox should be clicked successfully'

$ws.Cells.Item(25,3).Value = 'Enter valid email'
$ws.Cells.Item(25,4).Value = 'Email should be entered successfully(Entered email should be displayed in email inputbox)'

$ws.Cells.Item(26,3).Value = 'Click Phone Number'
$ws.Cells.Item(26,4).Value = 'Phon number inputbox should be clicked successfully'

$ws.Cells.Item(27,3).Value = 'Enter valid phone number'
$ws.Cells.Item(27,4).Value = 'Phon number should be enter successfully (Entered phon number should be displayed in Telephon inputbox'

$ws.Cells.Item(28,3).Value = 'Click password inputbox'
$ws.Cells.Item(28,4).Value = 'password inputbox should be clicked successfully'

$ws.Cells.Item(29,3).Value = 'Enter valid password'
$ws.Cells.Item(29,4).Value = 'Password should be enter successfully(some dots should be displayed in password inputbox)'

$ws.Cells.Item(30,3).Value = 'Click confirm password'
$ws.Cells.Item(30,4).Value = 'confirm password inputbox should be clicked successfully'

$ws.Cells.Item(31,3).Value = 'Enter confirm password (same as password)'
$ws.Cells.Item(31,4).Value = 'confirm password should be enter successfully (some dots should be displayed in confirm password inputbox)'

$ws.Cells.Item(32,3).Value = 'Click privacy Policy checkbox'
$ws.Cells.Item(32,4).Value = 'privacy policy checkbox should be displayed checked '

$ws.Cells.Item(33,3).Value = 'Click Continue button'
$ws.Cells.Item(33,4).Value = 'Error message should be displayed below blank mandatory filed (below Lastname inputbox because we should not fill lastname inputbox)'

# ---------------------------------------------------------------------------
# 3) TS_04 block (rows 34-48) - "Validate Register account fill all fields
#    but not select privacy policy Checkbox" - full flow but the privacy
#    policy checkbox is deliberately left unchecked.
# ---------------------------------------------------------------------------
$ws.Cells.Item(34,1).Value = 'TS_04'
$ws.Cells.Item(34,2).Value = 'Validate Register account fill all fields but not select privacy policy Checkbox'
$ws.Cells.Item(34,3).Value = 'Pre-requisite: TS_01'
$ws.Cells.Item(34,3).Font.Bold = $true
$ws.Cells.Item(34,3).Borders.LineStyle = 1
$ws.Cells.Item(34,3).HorizontalAlignment = -4131
$ws.Cells.Item(34,3).VerticalAlignment = -4160
$ws.Cells.Item(34,3).WrapText = $true

$ws.Cells.Item(35,3).Value = 'Click firstname inputbox'
$ws.Cells.Item(35,4).Value = 'firstname inputbox should be clicked successfully'

$ws.Cells.Item(36,3).Value = 'Enter valid firstname'
$ws.Cells.Item(36,4).Value = 'Firstname should be entered successfully (Entered firstname should be displayed in firstname inputbox)'

$ws.Cells.Item(37,3).Value = 'Click lastname inputbox'
$ws.Cells.Item(37,4).Value = 'lastname inputbox should be clicked successfully'

$ws.Cells.Item(38,3).Value = 'Enter valid lastname'
$ws.Cells.Item(38,4).Value = 'Lastname should be enterd successfully(Entered lastname should be displayed in lastname inputbox)'

$ws.Cells.Item(39,3).Value = 'Click email inputbox'
$ws.Cells.Item(39,4).Value = 'Email inputbox should be clicked successfully'

$ws.Cells.Item(40,3).Value = 'Enter valid email'
$ws.Cells.Item(40,4).Value = 'Email should be entered successfully(Entered email should be displayed in email inputbox)'

$ws.Cells.Item(41,3).Value = 'Click Phone Number'
$ws.Cells.Item(41,4).Value = 'Phon number inputbox should be clicked successfully'

$ws.Cells.Item(42,3).Value = 'Enter valid phone number'
$ws.Cells.Item(42,4).Value = 'Phon number should be enter successfully (Entered phon number should be displayed in Telephon inputbox'

$ws.Cells.Item(43,3).Value = 'Click password inputbox'
$ws.Cells.Item(43,4).Value = 'password inputbox should be clicked successfully'

$ws.Cells.Item(44,3).Value = 'Enter valid password'
$ws.Cells.Item(44,4).Value = 'Password should be enter successfully(some dots should be displayed in password inputbox)'

$ws.Cells.Item(45,3).Value = 'Click confirm password'
$ws.Cells.Item(45,4).Value = 'confirm password inputbox should be clicked successfully'

$ws.Cells.Item(46,3).Value = 'Enter confirm password (same as password)'
$ws.Cells.Item(46,4).Value = 'confirm password should be enter successfully (some dots should be displayed in confirm password inputbox)'

$ws.Cells.Item(47,3).Value = 'Not Click privacy Policy checkbox'
$ws.Cells.Item(47,4).Value = 'privacy policy checkbox should be displayed unchecked '

$ws.Cells.Item(48,3).Value = 'Click Continue button'
$ws.Cells.Item(48,4).Value = 'Your account has been created message should be displayed '

# ---------------------------------------------------------------------------
# 4) Row heights - row 10 got shorter now that column B is wider (less
#    wrapping), and every newly-filled row needs the height that its text
#    actually wraps to at the new column widths.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).RowHeight = 28.8

$ws.Rows.Item(20).RowHeight = 28.8
$ws.Rows.Item(21).RowHeight = 27
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(24).RowHeight = 25.2
$ws.Rows.Item(25).RowHeight = 28.8
$ws.Rows.Item(26).RowHeight = 28.8
$ws.Rows.Item(27).RowHeight = 43.2
$ws.Rows.Item(28).RowHeight = 23.4
$ws.Rows.Item(29).RowHeight = 28.8
$ws.Rows.Item(30).RowHeight = 28.8
$ws.Rows.Item(31).RowHeight = 43.2
$ws.Rows.Item(32).RowHeight = 28.8
$ws.Rows.Item(33).RowHeight = 43.2

$ws.Rows.Item(34).RowHeight = 28.8
$ws.Rows.Item(35).RowHeight = 27
$ws.Rows.Item(36).RowHeight = 43.2
$ws.Rows.Item(37).RowHeight = 24
$ws.Rows.Item(38).RowHeight = 43.2
$ws.Rows.Item(39).RowHeight = 25.2
$ws.Rows.Item(40).RowHeight = 28.8
$ws.Rows.Item(41).RowHeight = 28.8
$ws.Rows.Item(42).RowHeight = 43.2
$ws.Rows.Item(43).RowHeight = 23.4
$ws.Rows.Item(44).RowHeight = 28.8
$ws.Rows.Item(45).RowHeight = 28.8
$ws.Rows.Item(46).RowHeight = 43.2
$ws.Rows.Item(47).RowHeight = 28.8
$ws.Rows.Item(48).RowHeight = 28.8

# ---------------------------------------------------------------------------
# 5) The sheet no longer needs its very last blank row - remove it so the
#    one that used to be row 106 (the thicker bottom-border style) becomes
#    the new row 105, and keep the declared sheet dimension in sync.
# ---------------------------------------------------------------------------
$ws.Rows.Item(105).Delete()
$ws.Range("A105").Font.Bold = $ws.Range("A105").Font.Bold

# ---------------------------------------------------------------------------
# 6) View tidy-up: drop the frozen/scrolled-to top-left cell and move the
#    active selection to where the author was last working.
# ---------------------------------------------------------------------------
$ws.Range("C49").Select()
